$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 70

# Columns A and D hold text that Excel's type-inference would otherwise
# coerce into a date serial / a number ("2025-02-14" -> date, "06" -> 6).
# Force them to stay text by temporarily formatting as Text, then
# reset the style back to Normal (style index 0) so no stray number
# format sticks around on the cell.
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = "2025-02-14"
$ws.Cells.Item($r, 1).Style = "Normal"

$ws.Cells.Item($r, 2).Value = "22:05:18"
$ws.Cells.Item($r, 3).Value = "Friday"

$ws.Cells.Item($r, 4).NumberFormat = "@"
$ws.Cells.Item($r, 4).Value = "06"
$ws.Cells.Item($r, 4).Style = "Normal"

$ws.Cells.Item($r, 5).Value = 126029
$ws.Cells.Item($r, 6).Value = 142457
$ws.Cells.Item($r, 7).Value = 170518
$ws.Cells.Item($r, 8).Value = 159477
$ws.Cells.Item($r, 9).Value = -1
$ws.Cells.Item($r, 10).Value = 145087
$ws.Cells.Item($r, 11).Value = -1
$ws.Cells.Item($r, 12).Value = -1
$ws.Cells.Item($r, 13).Value = 192145
$ws.Cells.Item($r, 14).Value = 115098
$ws.Cells.Item($r, 15).Value = 45012
$ws.Cells.Item($r, 16).Value = 28707
$ws.Cells.Item($r, 17).Value = 66012
$ws.Cells.Item($r, 18).Value = -1
$ws.Cells.Item($r, 19).Value = 46161
$ws.Cells.Item($r, 20).Value = -1
